$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the previous data row down into the two new rows
$ws.Range("A46:C46").Copy()
$ws.Range("A48:C49").PasteSpecial(-4122)

# Fill in the two new localization rows (Key, String EN, String DE)
$ws.Range("A48").Value = "service_worker-update_headline"
$ws.Range("A49").Value = "service_worker-update_confirm_btn_txt"
$ws.Range("B48").Value = "Update Available"
$ws.Range("C48").Value = "Update verfügbar"
$ws.Range("B49").Value = "Update Now & Refresh"
$ws.Range("C49").Value = "Update installieren"

# Grow the "Tabelle2" table so it covers the two new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C49"))

# Re-fit column A now that it holds slightly longer key strings
$ws.Columns.Item(1).ColumnWidth = 31.8

# Move the viewport / selection like the saved workbook shows
$null = $ws.Range("D44").Select()
